$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2:G5").Value = 0.057308
$ws.Range("H2:H5").Value = 0.171924

$ws.Range("M2").Value = 29.04796866666666
$ws.Range("N2").Value = 87.14390599999999

$ws.Range("M4").Value = 12.70280433333333
$ws.Range("N4").Value = 38.108413

$ws.Range("M5").Value = 48.89716966666666
$ws.Range("N5").Value = 146.691509

$ws.Range("O2").Value = 0.2371972210028098
$ws.Range("P2").Value = 0.2371972210028099
$ws.Range("Q2").Value = 1.664680988349333
$ws.Range("R2").Value = 14.982128895144
$ws.Range("S2").Value = 0.2371972210028098
$ws.Range("T2").Value = 0.2371972210028099

$ws.Range("O3").Value = 0.2597953978506987
$ws.Range("P3").Value = 0.2597953978506987
$ws.Range("Q3").Value = 1.823277936538667
$ws.Range("R3").Value = 16.409501428848
$ws.Range("S3").Value = 0.2597953978506987
$ws.Range("T3").Value = 0.2597953978506987

$ws.Range("O4").Value = 0.1037273869778955
$ws.Range("P4").Value = 0.1037273869778955
$ws.Range("Q4").Value = 0.7279723107346666
$ws.Range("R4").Value = 6.551750796612
$ws.Range("S4").Value = 0.1037273869778955
$ws.Range("T4").Value = 0.1037273869778955

$ws.Range("O5").Value = 0.3992799941685959
$ws.Range("P5").Value = 0.399279994168596
$ws.Range("Q5").Value = 2.802198999257333
$ws.Range("R5").Value = 25.219790993316
$ws.Range("S5").Value = 0.3992799941685959
$ws.Range("T5").Value = 0.399279994168596
